$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 165; this pushes the previous
# rows 165-180 down to become rows 167-182 (matching the target dimension
# A1:R182).
$ws.Rows(165).Insert()
$ws.Rows(165).Insert()

# --- New row 165 (weekly "Americana (o)" entry, Región Metropolitana) ---
$ws.Range("A165").Value = 11
$ws.Range("B165").Value = "Vega Monumental Concepción"
$ws.Range("C165").Value = "Bíobío"
$ws.Range("D165").Value = 45013
$ws.Range("E165").Value = 8
$ws.Range("F165").Value = 100112021
$ws.Range("G165").Value = "Ají"
$ws.Range("H165").Value = "Americana (o)"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 110
$ws.Range("K165").Value = 18000
$ws.Range("L165").Value = 19000
$ws.Range("M165").Value = 18455
$ws.Range("N165").Value = "$/saco 25 kilos"
$ws.Range("O165").Value = "Región Metropolitana"
$ws.Range("P165").Value = 738
$ws.Range("Q165").Value = 25
$ws.Range("R165").Value = "Hortaliza"

# --- New row 166 (weekly "Chilena(o)" entry, Región Metropolitana) ---
$ws.Range("A166").Value = 11
$ws.Range("B166").Value = "Vega Monumental Concepción"
$ws.Range("C166").Value = "Bíobío"
$ws.Range("D166").Value = 45013
$ws.Range("E166").Value = 8
$ws.Range("F166").Value = 100112021
$ws.Range("G166").Value = "Ají"
$ws.Range("H166").Value = "Chilena(o)"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 30
$ws.Range("K166").Value = 13000
$ws.Range("L166").Value = 14000
$ws.Range("M166").Value = 13500
$ws.Range("N166").Value = "$/caja 15 kilos"
$ws.Range("O166").Value = "Región Metropolitana"
$ws.Range("P166").Value = 900
$ws.Range("Q166").Value = 15
$ws.Range("R166").Value = "Hortaliza"
